# Weekly update: a new price record is inserted into the "Ají" (Hortaliza)
# sheet. The new record is inserted as row 35, pushing the previously
# existing rows 35-83 down to rows 36-84 (dimension grows from A1:R83 to
# A1:R84).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 35; Excel shifts row 35
# and everything below it down by one, preserving their values/styles.
$ws.Rows.Item(35).Insert()

# Fill in the data for the newly inserted row 35.
$ws.Range("A35").Value = 1
$ws.Range("B35").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C35").Value = "Arica y Parinacota"
$ws.Range("D35").Value = 44725
$ws.Range("E35").Value = 15
$ws.Range("F35").Value = 100112021
$ws.Range("G35").Value = "Ají"
$ws.Range("H35").Value = "Inferno"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 150
$ws.Range("K35").Value = 18000
$ws.Range("L35").Value = 19000
$ws.Range("M35").Value = 18500
$ws.Range("N35").Value = "$/caja 15 kilos"
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value = 1233
$ws.Range("Q35").Value = 15
$ws.Range("R35").Value = "Hortaliza"
